# "correzione formula nastro trasportatore"
# Fix the conveyor-belt ("Nastro Trasportatore") sizing formulas on the
# "Design1" sheet:
#   - the belt-friction-table selector (C12) was left on the wrong row
#     factor (1/110 instead of 1/55) -> switch it from 1 to 2, which
#     ripples through the whole resistance/power/motor computation below.
#   - the r3 resistance (C29) was referencing the wrong coefficient row
#     (C16) instead of the one actually driven by the belt width (C19).
#   - the number-of-rollers rounding (C44) hard-coded a "5" instead of
#     referencing the roller-spacing input cell (C43).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Design1")

$ws.Range("C12").Value = 2
$ws.Range("C29").Formula = "=C19*C6*C5*9.81"
$ws.Range("C44").Formula = "=_xlfn.CEILING.MATH(C38/(C43*B19))"

# Matches the cursor position left behind in the saved workbook.
[void]$ws.Range("G57").Select()
